$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New player data (rows 2-16 reshuffled; rows 17-19 unchanged)
$data = @(
    @("Dyson Daniels", "PG,SG", "Atlanta Hawks"),
    @("Malik Beasley", "SG", "Detroit Pistons"),
    @("Jamal Murray", "PG,SG", "Denver Nuggets"),
    @("Ochai Agbaji", "SG,SF", "Toronto Raptors"),
    @("Josh Hart", "SF,PF", "New York Knicks"),
    @("De'Andre Hunter", "SF,PF", "Atlanta Hawks"),
    @("Tari Eason", "SF,PF", "Houston Rockets"),
    @("Kyle Filipowski", "PF,C", "Utah Jazz"),
    @("Alexandre Sarr", "PF,C", "Washington Wizards"),
    @("Domantas Sabonis", "C", "Sacramento Kings"),
    @("Myles Turner", "C", "Indiana Pacers"),
    @("Victor Wembanyama", "C", "San Antonio Spurs"),
    @("Luguentz Dort", "SG,SF", "Oklahoma City Thunder"),
    @("Naz Reid", "PF,C", "Minnesota Timberwolves"),
    @("Donovan Mitchell", "PG,SG", "Cleveland Cavaliers")
)

$row = 2
foreach ($player in $data) {
    $ws.Cells.Item($row, 1).Value = $player[0]
    $ws.Cells.Item($row, 2).Value = $player[1]
    $ws.Cells.Item($row, 3).Value = $player[2]
    $row++
}
